# Adds 17 new departure rows (rows 499-515) to the "Main Data" sheet,
# mirroring the new flights appended to the KRK_Departures table for
# Saturday, Jan 14 (columns A-J and L; K and M stay blank as in the
# existing rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: Number, Date, Time, Flight, To, Short, Airline, Model, AircraftId, Status, Difference
$newDepartures = @(
  @(498.0, "Saturday, Jan 14", "5:45 PM", "W65003", "London", "(LTN)", "Wizz Air ", "A321", "(HA-LXO)", "6:01 PM", "0 hours, 16 minutes"),
  @(499.0, "Saturday, Jan 14", "5:50 PM", "FR4935", "Brussels", "(CRL)", "Ryanair ", "B738", "(SP-RKC)", "6:06 PM", "0 hours, 16 minutes"),
  @(500.0, "Saturday, Jan 14", "5:55 PM", "FR1812", "London", "(LTN)", "Ryanair ", "B38M", "(SP-RZK)", "5:57 PM", "0 hours, 2 minutes"),
  @(501.0, "Saturday, Jan 14", "6:25 PM", "FR6240", "Naples", "(NAP)", "Buzz ", "B38M", "(SP-RZF)", "6:58 PM", "0 hours, 33 minutes"),
  @(502.0, "Saturday, Jan 14", "6:40 PM", "LH1369", "Frankfurt", "(FRA)", "Eurowings ", "A319", "(D-AGWL)", "6:48 PM", "0 hours, 8 minutes"),
  @(503.0, "Saturday, Jan 14", "6:55 PM", "BA873", "London", "(LHR)", "British Airways ", "A320", "(G-MIDO)", "7:14 PM", "0 hours, 19 minutes"),
  @(504.0, "Saturday, Jan 14", "6:55 PM", "FR6216", "Oslo", "(TRF)", "Buzz ", "B38M", "(SP-RZD)", "7:16 PM", "0 hours, 21 minutes"),
  @(505.0, "Saturday, Jan 14", "6:55 PM", "W65061", "Rome", "(FCO)", "Wizz Air ", "A21N", "(HA-LZI)", "7:02 PM", "0 hours, 7 minutes"),
  @(506.0, "Saturday, Jan 14", "7:00 PM", "FR6244", "Lille", "(LIL)", "Buzz ", "B38M", "(SP-RZC)", "7:40 PM", "0 hours, 40 minutes"),
  @(507.0, "Saturday, Jan 14", "7:05 PM", "FR1643", "Vienna", "(VIE)", "Ryanair ", "B738", "(SP-RST)", "7:24 PM", "0 hours, 19 minutes"),
  @(508.0, "Saturday, Jan 14", "7:10 PM", "FR6252", "Stockholm", "(ARN)", "Ryanair ", "B738", "(SP-RKB)", "7:32 PM", "0 hours, 22 minutes"),
  @(509.0, "Saturday, Jan 14", "7:50 PM", "FR6276", "Nuremberg", "(NUE)", "Ryanair ", "B738", "(9H-QCX)", "7:49 PM", "0 hours, -1 minutes"),
  @(510.0, "Saturday, Jan 14", "8:50 PM", "FR7969", "Pisa", "(PSA)", "Ryanair ", "B738", "(9H-QDG)", "9:02 PM", "0 hours, 12 minutes"),
  @(511.0, "Saturday, Jan 14", "9:05 PM", "FR6230", "Szczecin", "(SZZ)", "Buzz ", "B38M", "(SP-RZB)", "9:29 PM", "0 hours, 24 minutes"),
  @(512.0, "Saturday, Jan 14", "9:15 PM", "FR2333", "Leeds", "(LBA)", "Ryanair ", "B738", "(EI-ENL)", "9:34 PM", "0 hours, 19 minutes"),
  @(513.0, "Saturday, Jan 14", "9:15 PM", "FR2713", "London", "(STN)", "Ryanair ", "B738", "(EI-EKN)", "9:38 PM", "0 hours, 23 minutes"),
  @(514.0, "Saturday, Jan 14", "9:50 PM", "FR6248", "Manchester", "(MAN)", "Ryanair ", "B738", "(EI-DHB)", "10:02 PM", "0 hours, 12 minutes")
)

$row = 499
foreach ($flight in $newDepartures) {
  $ws.Cells.Item($row, 1).Value = $flight[0]   # A: NUMBER
  $ws.Cells.Item($row, 2).Value = $flight[1]   # B: DATE
  $ws.Cells.Item($row, 3).Value = $flight[2]   # C: TIME
  $ws.Cells.Item($row, 4).Value = $flight[3]   # D: FLIGHT
  $ws.Cells.Item($row, 5).Value = $flight[4]   # E: TO
  $ws.Cells.Item($row, 6).Value = $flight[5]   # F: SHORT
  $ws.Cells.Item($row, 7).Value = $flight[6]   # G: AIRLINE
  $ws.Cells.Item($row, 8).Value = $flight[7]   # H: MODEL
  $ws.Cells.Item($row, 9).Value = $flight[8]   # I: AIRCFAT ID
  $ws.Cells.Item($row, 10).Value = $flight[9]  # J: STATUS
  $ws.Cells.Item($row, 12).Value = $flight[10] # L: DIFFERENCE (K and M left blank)
  $row = $row + 1
}
